$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: new product name and amount
$ws.Range("A2").Value = "Fastrack Stunners 1.0"
$ws.Range("B2").Value = 12340.5

# Remove the old row 3 (Fastrack Streetline 3.0 Analog-Digital Watch) entirely,
# shifting the "Final Total" row up to become row 3.
$ws.Rows(3).Delete()

# Update the Final Total amount to match the new total.
$ws.Range("B3").Value = 12340.5
